$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 11959.0205
$ws.Range("I15").Value = 11959.0205
$ws.Range("K15").Value = 35877.0615
$ws.Range("M15").Value = -35708.0615

$ws.Range("H17").Value = 5013686
$ws.Range("J17").Value = 5013686
$ws.Range("L17").Value = 15041058
$ws.Range("N17").Value = -15041394

$ws.Range("H132").Value = 12822524
$ws.Range("I132").Value = 1444.5333
$ws.Range("J132").Value = 30305814
$ws.Range("K132").Value = 4333.5999
$ws.Range("L132").Value = 90917442
$ws.Range("M132").Value = -1803.5999
$ws.Range("N132").Value = -90922502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8593.615
$ws.Range("I37").Value = 7545.4287
$ws.Range("J37").Value = 9816.5
$ws.Range("K37").Value = 7545.4287
$ws.Range("L37").Value = 9816.5
$ws.Range("M37").Value = -7272.4287
$ws.Range("N37").Value = -10362.5

$ws.Range("H44").Value = 18066.666
$ws.Range("J44").Value = 18066.666
$ws.Range("L44").Value = 18066.666
$ws.Range("N44").Value = -19042.666

$ws.Range("H55").Value = 20075
$ws.Range("J55").Value = 20075
$ws.Range("L55").Value = 20075
$ws.Range("N55").Value = -20705

$ws.Range("H61").Value = 3179.6948
$ws.Range("I61").Value = 3362.5293
$ws.Range("J61").Value = 2014.125
$ws.Range("K61").Value = 3362.5293
$ws.Range("L61").Value = 2014.125
$ws.Range("M61").Value = -3150.5293
$ws.Range("N61").Value = -2438.125

$ws.Range("H80").Value = 24943.5
$ws.Range("I80").Value = 7777
$ws.Range("K80").Value = 7777
$ws.Range("M80").Value = -6779

$ws.Range("H83").Value = 24943.5
$ws.Range("I83").Value = 7777
$ws.Range("K83").Value = 23331
$ws.Range("M83").Value = -18339

$ws.Range("H122").Value = 1427578.5
$ws.Range("I122").Value = 1605694.5
$ws.Range("J122").Value = 2650
$ws.Range("K122").Value = 4817083.5
$ws.Range("L122").Value = 7950
$ws.Range("M122").Value = -4814633.5
$ws.Range("N122").Value = -12850

$ws.Range("H132").Value = 2312.8372
$ws.Range("I132").Value = 1387.5
$ws.Range("J132").Value = 5808.5557
$ws.Range("K132").Value = 4162.5
$ws.Range("L132").Value = 17425.6671
$ws.Range("M132").Value = -1632.5
$ws.Range("N132").Value = -22485.6671

$ws.Range("H136").Value = 3179.6948
$ws.Range("I136").Value = 3362.5293
$ws.Range("J136").Value = 2014.125
$ws.Range("K136").Value = 10087.5879
$ws.Range("L136").Value = 6042.375
$ws.Range("M136").Value = -7537.5879
$ws.Range("N136").Value = -11142.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3841.68
$ws.Range("I134").Value = 4892.484
$ws.Range("K134").Value = 14677.452
$ws.Range("M134").Value = -12142.452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13164740
$ws.Range("I31").Value = 2755.375
$ws.Range("J31").Value = 22737092
$ws.Range("K31").Value = 2755.375
$ws.Range("L31").Value = 22737092
$ws.Range("M31").Value = -2460.375
$ws.Range("N31").Value = -22737682

$ws.Range("H34").Value = 13164740
$ws.Range("I34").Value = 2755.375
$ws.Range("J34").Value = 22737092
$ws.Range("K34").Value = 2755.375
$ws.Range("L34").Value = 22737092
$ws.Range("M34").Value = -2553.375
$ws.Range("N34").Value = -22737496

$ws.Range("H58").Value = 5210134.5
$ws.Range("I58").Value = 9260425
$ws.Range("J58").Value = 2617.2856
$ws.Range("K58").Value = 9260425
$ws.Range("L58").Value = 2617.2856
$ws.Range("M58").Value = -9260222
$ws.Range("N58").Value = -3023.2856

$ws.Range("H132").Value = 6899343.5
$ws.Range("I132").Value = 9525951
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 28577853
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -28575323
$ws.Range("N132").Value = -18555.5

$ws.Range("H136").Value = 5210134.5
$ws.Range("I136").Value = 9260425
$ws.Range("J136").Value = 2617.2856
$ws.Range("K136").Value = 27781275
$ws.Range("L136").Value = 7851.8568
$ws.Range("M136").Value = -27778725
$ws.Range("N136").Value = -12951.8568

$ws.Range("H140").Value = 23034.281
$ws.Range("J140").Value = 23034.281
$ws.Range("L140").Value = 23034.281
$ws.Range("N140").Value = -33394.281

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 80
$ws.Range("I50").Value = 80
$ws.Range("K50").Value = 240
$ws.Range("M50").Value = 241

$ws.Range("H53").Value = 80
$ws.Range("I53").Value = 80
$ws.Range("K53").Value = 240
$ws.Range("M53").Value = 241

$ws.Range("H117").Value = 15878826
$ws.Range("I117").Value = 363.6
$ws.Range("J117").Value = 20840846
$ws.Range("K117").Value = 1090.8
$ws.Range("L117").Value = 62522538
$ws.Range("M117").Value = 2351.2
$ws.Range("N117").Value = -62529422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7248933.5
$ws.Range("I132").Value = 10418942
$ws.Range("J132").Value = 3198.2856
$ws.Range("K132").Value = 31256826
$ws.Range("L132").Value = 9594.856800000001
$ws.Range("M132").Value = -31254296
$ws.Range("N132").Value = -14654.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4050
$ws.Range("I61").Value = 3175
$ws.Range("J61").Value = 4925
$ws.Range("K61").Value = 3175
$ws.Range("L61").Value = 4925
$ws.Range("M61").Value = -2973
$ws.Range("N61").Value = -5329

$ws.Range("H113").Value = 4050
$ws.Range("I113").Value = 3175
$ws.Range("J113").Value = 4925
$ws.Range("K113").Value = 3175
$ws.Range("L113").Value = 4925
$ws.Range("M113").Value = -1005
$ws.Range("N113").Value = -9265

$ws.Range("H132").Value = 9756626
$ws.Range("I132").Value = 12066393
$ws.Range("J132").Value = 4278.1113
$ws.Range("K132").Value = 36199179
$ws.Range("L132").Value = 12834.3339
$ws.Range("M132").Value = -36196649
$ws.Range("N132").Value = -17894.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7013.636
$ws.Range("I54").Value = 2950
$ws.Range("K54").Value = 2950
$ws.Range("M54").Value = -2430

$ws.Range("H81").Value = 40001560
$ws.Range("J81").Value = 100002500
$ws.Range("L81").Value = 200005000
$ws.Range("N81").Value = -200007122

$ws.Range("H84").Value = 40001560
$ws.Range("J84").Value = 100002500
$ws.Range("L84").Value = 1000025000
$ws.Range("N84").Value = -1000035608
